# Updates quest-leve market price / profit figures across all job sheets
# (data refresh from scheduled market-board pull).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 83334690
$ws.Range("I12").Value = 1328.1
$ws.Range("J12").Value = 500001500
$ws.Range("K12").Value = 1328.1
$ws.Range("L12").Value = 500001500
$ws.Range("M12").Value = -1158.1
$ws.Range("N12").Value = -500001840
$ws.Range("H74").Value = 4324.6875
$ws.Range("I74").Value = 4169
$ws.Range("K74").Value = 4169
$ws.Range("M74").Value = -3233
$ws.Range("H77").Value = 4324.6875
$ws.Range("I77").Value = 4169
$ws.Range("K77").Value = 20845
$ws.Range("M77").Value = -16165
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H130").Value = 49980
$ws.Range("J130").Value = 49980
$ws.Range("L130").Value = 49980
$ws.Range("N130").Value = -60020

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 101567.3
$ws.Range("I74").Value = 126320.125
$ws.Range("J74").Value = 2556
$ws.Range("K74").Value = 126320.125
$ws.Range("L74").Value = 2556
$ws.Range("M74").Value = -125446.125
$ws.Range("N74").Value = -4304
$ws.Range("H77").Value = 101567.3
$ws.Range("I77").Value = 126320.125
$ws.Range("J77").Value = 2556
$ws.Range("K77").Value = 631600.625
$ws.Range("L77").Value = 12780
$ws.Range("M77").Value = -627232.625
$ws.Range("N77").Value = -21516
$ws.Range("H92").Value = 28750
$ws.Range("J92").Value = 28750
$ws.Range("L92").Value = 28750
$ws.Range("N92").Value = -33742

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1803.1
$ws.Range("I86").Value = 1603.9286
$ws.Range("J86").Value = 2267.8333
$ws.Range("K86").Value = 1603.9286
$ws.Range("L86").Value = 2267.8333
$ws.Range("M86").Value = -480.9286
$ws.Range("N86").Value = -4513.8333
$ws.Range("H89").Value = 1803.1
$ws.Range("I89").Value = 1603.9286
$ws.Range("J89").Value = 2267.8333
$ws.Range("K89").Value = 8019.643
$ws.Range("L89").Value = 11339.1665
$ws.Range("M89").Value = -2403.643
$ws.Range("N89").Value = -22571.1665
$ws.Range("H134").Value = 5148.7144
$ws.Range("I134").Value = 4936.9565
$ws.Range("J134").Value = 6122.8
$ws.Range("K134").Value = 14810.8695
$ws.Range("L134").Value = 18368.4
$ws.Range("M134").Value = -12275.8695
$ws.Range("N134").Value = -23438.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 7369.727
$ws.Range("I32").Value = 1940
$ws.Range("J32").Value = 10472.429
$ws.Range("K32").Value = 1940
$ws.Range("L32").Value = 10472.429
$ws.Range("M32").Value = -1624
$ws.Range("N32").Value = -11104.429
$ws.Range("H45").Value = 6445
$ws.Range("I45").Value = 3000
$ws.Range("J45").Value = 9890
$ws.Range("K45").Value = 3000
$ws.Range("L45").Value = 9890
$ws.Range("M45").Value = -2407
$ws.Range("N45").Value = -11076

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 825.9
$ws.Range("I6").Value = 656.875
$ws.Range("J6").Value = 1502
$ws.Range("K6").Value = 1970.625
$ws.Range("L6").Value = 4506
$ws.Range("M6").Value = -1857.625
$ws.Range("N6").Value = -4732
$ws.Range("H70").Value = 6349
$ws.Range("I70").Value = 5610
$ws.Range("J70").Value = 6771.2856
$ws.Range("K70").Value = 16830
$ws.Range("L70").Value = 20313.8568
$ws.Range("M70").Value = -16515
$ws.Range("N70").Value = -20943.8568
$ws.Range("H73").Value = 6349
$ws.Range("I73").Value = 5610
$ws.Range("J73").Value = 6771.2856
$ws.Range("K73").Value = 16830
$ws.Range("L73").Value = 20313.8568
$ws.Range("M73").Value = -15738
$ws.Range("N73").Value = -22497.8568
$ws.Range("H92").Value = 1667005
$ws.Range("I92").Value = 420
$ws.Range("J92").Value = 5000175
$ws.Range("K92").Value = 1260
$ws.Range("L92").Value = 15000525
$ws.Range("M92").Value = -12
$ws.Range("N92").Value = -15003021
$ws.Range("H99").Value = 3510
$ws.Range("I99").Value = 2775
$ws.Range("K99").Value = 8325
$ws.Range("M99").Value = -6079
$ws.Range("H103").Value = 3168.149
$ws.Range("I103").Value = 459.16666
$ws.Range("J103").Value = 3564.5854
$ws.Range("K103").Value = 1377.49998
$ws.Range("L103").Value = 10693.7562
$ws.Range("M103").Value = -498.4999800000001
$ws.Range("N103").Value = -12451.7562

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 20.38889
$ws.Range("I2").Value = 10.357142
$ws.Range("J2").Value = 55.5
$ws.Range("K2").Value = 10.357142
$ws.Range("L2").Value = 55.5
$ws.Range("M2").Value = 102.642858
$ws.Range("N2").Value = -281.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2057.2856
$ws.Range("I82").Value = 850
$ws.Range("J82").Value = 2341.353
$ws.Range("K82").Value = 850
$ws.Range("L82").Value = 2341.353
$ws.Range("M82").Value = -489
$ws.Range("N82").Value = -3063.353
$ws.Range("H85").Value = 2057.2856
$ws.Range("I85").Value = 850
$ws.Range("J85").Value = 2341.353
$ws.Range("K85").Value = 850
$ws.Range("L85").Value = 2341.353
$ws.Range("M85").Value = 398
$ws.Range("N85").Value = -4837.353
$ws.Range("H93").Value = 34677.906
$ws.Range("I93").Value = 1454.8235
$ws.Range("J93").Value = 175876
$ws.Range("K93").Value = 1454.8235
$ws.Range("L93").Value = 175876
$ws.Range("M93").Value = -206.8235
$ws.Range("N93").Value = -178372
$ws.Range("H127").Value = 33606.4
$ws.Range("J127").Value = 33606.4
$ws.Range("L127").Value = 33606.4
$ws.Range("N127").Value = -43526.4
$ws.Range("H128").Value = 33999.855
$ws.Range("J128").Value = 33999.855
$ws.Range("L128").Value = 33999.855
$ws.Range("N128").Value = -43959.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 27000
$ws.Range("J64").Value = 27000
$ws.Range("L64").Value = 27000
$ws.Range("N64").Value = -27496
$ws.Range("H67").Value = 27000
$ws.Range("J67").Value = 27000
$ws.Range("L67").Value = 27000
$ws.Range("N67").Value = -28716
$ws.Range("H135").Value = 38000
$ws.Range("J135").Value = 38000
$ws.Range("L135").Value = 38000
$ws.Range("N135").Value = -48140
